$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '66.502.65'
$ws.Range('E2').Value = '  -3.66%  '
$ws.Range('D3').Value = '3.562.17'
$ws.Range('E3').Value = '  -4.34%  '
$ws.Range('E4').Value = '  +0.10%  '
Set-TextValue 'D5' '572.37'
$ws.Range('E5').Value = '  -6.92%  '
Set-TextValue 'D6' '186.22'
$ws.Range('E6').Value = '  -3.57%  '
$ws.Range('D7').Value = '3.560.13'
$ws.Range('E7').Value = '  -4.23%  '
$ws.Range('E8').Value = '  -4.10%  '
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('E10').Value = '  -7.48%  '
Set-TextValue 'D11' '55.93'
$ws.Range('E11').Value = '  -7.28%  '
$ws.Range('E12').Value = '  -7.83%  '
$ws.Range('E13').Value = '  -10.08%  '
$ws.Range('E14').Value = '  -6.35%  '
$ws.Range('D15').Value = '4.129.42'
$ws.Range('E15').Value = '  -4.29%  '
$ws.Range('D16').Value = '3.569.11'
$ws.Range('E16').Value = '  -4.17%  '
$ws.Range('E17').Value = '  -1.77%  '
Set-TextValue 'D18' '18.26'
$ws.Range('E18').Value = '  -6.46%  '
$ws.Range('D19').Value = '66.526.84'
$ws.Range('E19').Value = '  -3.43%  '
Set-TextValue 'D20' '12.08'
$ws.Range('E20').Value = '  -7.00%  '
$ws.Range('E21').Value = '  -8.56%  '
Set-TextValue 'D22' '388.93'
$ws.Range('E22').Value = '  -5.90%  '
Set-TextValue 'D23' '4.19'
$ws.Range('E23').Value = '  -8.64%  '
Set-TextValue 'D24' '85.38'
$ws.Range('E24').Value = '  -5.45%  '
Set-TextValue 'D25' '11.22'
$ws.Range('E25').Value = '  -1.46%  '
Set-TextValue 'D26' '2.92'
$ws.Range('E26').Value = '  -6.40%  '
$ws.Range('E27').Value = '  -6.06%  '
$ws.Range('E28').Value = '  +0.00%  '
Set-TextValue 'D29' '3.55'
$ws.Range('E29').Value = '  -7.47%  '
Set-TextValue 'D30' '8.82'
$ws.Range('E30').Value = '  -9.51%  '
Set-TextValue 'D31' '7.55'
$ws.Range('E31').Value = '  -3.22%  '
Set-TextValue 'D32' '30.81'
$ws.Range('E32').Value = '  -6.60%  '
Set-TextValue 'D33' '626.79'
$ws.Range('E33').Value = '  -2.16%  '
Set-TextValue 'D34' '12.15'
$ws.Range('E34').Value = '  -4.97%  '
$ws.Range('E35').Value = '  -7.88%  '
Set-TextValue 'D36' '63.30'
$ws.Range('E36').Value = '  -6.36%  '
Set-TextValue 'D37' '41.78'
$ws.Range('E37').Value = '  -10.87%  '
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D38' '0.403'
$ws.Range('E38').Value = '  -3.39%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D39' '1.00'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').Value = '0.0₃0746'
$ws.Range('E40').Value = '  -10.64%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.141.68'
$ws.Range('E41').Value = '  +7.17%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D42' '0.133'
$ws.Range('E42').Value = '  -5.56%  '
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('E44').Value = '  -4.33%  '
Set-TextValue 'D45' '2.65'
$ws.Range('E45').Value = '  +0.49%  '
Set-TextValue 'D46' '0.0411'
$ws.Range('E46').Value = '  -8.61%  '
$ws.Range('E47').Value = '  -7.10%  '
Set-TextValue 'D48' '3.04'
$ws.Range('E48').Value = '  -2.09%  '
Set-TextValue 'D49' '139.18'
$ws.Range('E49').Value = '  -4.02%  '
Set-TextValue 'D50' '8.41'
$ws.Range('E50').Value = '  -9.97%  '
$ws.Range('E51').Value = '  -2.08%  '
